$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano test statistics (C) and p-values (D)
$ws.Range("C2").Value = -1.419589570654476
$ws.Range("D2").Value = 0.1697424384815207

$ws.Range("C3").Value = 0.2018720937891392
$ws.Range("D3").Value = 0.8418717398055959

$ws.Range("C4").Value = 0.1479249450853262
$ws.Range("D4").Value = 0.8837494078178698

$ws.Range("C5").Value = -0.8936201948015642
$ws.Range("D5").Value = 0.3811970180824504

$ws.Range("C6").Value = 1.486649141724645
$ws.Range("D6").Value = 0.1512994932407943

$ws.Range("C7").Value = 1.885000145528612
$ws.Range("D7").Value = 0.07270783678968828

$ws.Range("C8").Value = 0.9103759789715122
$ws.Range("D8").Value = 0.3724913412073843

$ws.Range("C9").Value = -0.02949473566067731
$ws.Range("D9").Value = 0.9767359525653478

$ws.Range("C10").Value = -0.9557333443835314
$ws.Range("D10").Value = 0.3495919743532325

$ws.Range("C11").Value = -0.8976064817577675
$ws.Range("D11").Value = 0.3791138700936876
